$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $rng = $d.Paragraphs($index).Range
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-InParagraph 62 "Gorofa" "00:01:28,000 --> 00:01:32,000"
Replace-InParagraph 63 "peak is one meter wide the two ants move" "kilele ni mita moja upana wa mchwa wawili hoja"
Replace-InParagraph 91 "juu ya mlima. Your purpose is to" "juu ya mlima. Kusudi lako ni"
Replace-InParagraph 95 "make the time the last ant takes before" "fanya wakati mchwa wa mwisho huchukua hapo awali"
Replace-InParagraph 99 "falling the longest possible. Ants cannot" "kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi"
Replace-InParagraph 103 "be still: they must move to the right or" "tulia: lazima wahamie kulia au"
Replace-InParagraph 107 "to the left but they must move and after" "upande wa kushoto lakini lazima wasogee na baada"
Replace-InParagraph 111 "meeting each other they turn around and" "wakikutana wanageuka na"
Replace-InParagraph 115 "keep moving with the same but opposite" "endelea kusonga na sawa lakini kinyume"
Replace-InParagraph 119 "velocity" "kasi"
Replace-InParagraph 127 "so again what are the precise positions" "kwa hivyo tena ni nafasi gani sahihi"
Replace-InParagraph 131 "where I should place the two ants in" "ambapo ninapaswa kuwaweka mchwa wawili ndani"
Replace-InParagraph 135 "order to get the longest time before the" "ili kupata muda mrefu zaidi kabla ya"
Replace-InParagraph 139 "last ant falls? The second puzzle is" "chungu mwisho huanguka? Fumbo la pili ni"
Replace-InParagraph 143 "basically the same but now we have three" "kimsingi ni sawa lakini sasa tuna tatu"
Replace-InParagraph 147 "ants instead of two." "mchwa badala ya wawili."
Replace-InParagraph 151 "As before the ants velocity is one" "Kama kabla ya mchwa kasi ni moja"
Replace-InParagraph 155 "centimeter per second, every ant turns" "sentimita kwa sekunde, kila mchwa hugeuka"
Replace-InParagraph 159 "around after meeting another ant and" "karibu baada ya kukutana na mchwa mwingine na"
Replace-InParagraph 163 "the peak is one meter wide. So, what are" "kilele kina upana wa mita moja. Hivyo, ni nini"
Replace-InParagraph 167 "now the precise positions" "sasa nafasi sahihi"
Replace-InParagraph 171 "I should place the three ants in order" "Ninapaswa kuweka mchwa watatu kwa mpangilio"
Replace-InParagraph 175 "to get the longest time before the last" "kupata muda mrefu zaidi kabla ya mwisho"
Replace-InParagraph 179 "ant falls down? I hope you enjoyed this" "chungu huanguka chini? Natumaini ulifurahia hili"
Replace-InParagraph 183 "video do your best and good luck" "video fanya bora na bahati nzuri"
